$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 7-15 with new GR numbers
$ws.Range("A7").Value = 21810496
$ws.Range("A8").Value = 21810685
$ws.Range("A9").Value = 21810293
$ws.Range("A10").Value = 21920031
$ws.Range("A11").Value = 21810244
$ws.Range("A12").Value = 21810277
$ws.Range("A13").Value = 21810275
$ws.Range("A14").Value = 21810417
$ws.Range("A15").Value = 21810282

# Add new rows 16-37 with new GR numbers, copying style from A15
$newValues = @(21810878, 21920063, 21810308, 21920157, 21810099, 21810149, 21810427, 21810824, 21810251, 21810565, 21810827, 21810169, 21810535, 21810715, 21810483, 21810228, 21810364, 21810339, 21810590, 21810062, 21810261, 21810703)

$row = 16
foreach ($val in $newValues) {
    $cell = $ws.Cells.Item($row, 1)
    $ws.Range("A15").Copy()
    $cell.PasteSpecial(-4122)  # xlPasteFormats
    $cell.Value = $val
    $row = $row + 1
}
$excel.CutCopyMode = $false

# Update selection / view to match target (scrolled so row 31 is the top visible row)
$ws.Range("B36").Select()
$excel.ActiveWindow.ScrollRow = 31
$excel.ActiveWindow.ScrollColumn = 1
